$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 639 ("「彼は応えるだろう」...") was removed from the post list.
# Deleting the entire row shifts all subsequent rows (640-807) up by one,
# which matches the renumbering seen across the rest of the sheet, and
# Excel will automatically adjust the sheet dimension to A1:C806.
$ws.Rows("639").Delete()
